$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows starting at row 138, pushing the existing
# 20100203.. data (and everything below it) down by 4 rows.
$ws.Rows.Item(138).Resize(4).Insert()

# Make the new date cells behave like the rest of column A (plain text
# dates), matching the existing data's cell type.
$ws.Range("A138:A141").NumberFormat = "@"

$ws.Cells.Item(138,1).Value = "20091222"
$ws.Cells.Item(138,2).Value = 776

$ws.Cells.Item(139,1).Value = "20100104"
$ws.Cells.Item(139,2).Value = 856

$ws.Cells.Item(140,1).Value = "20100108"
$ws.Cells.Item(140,2).Value = 787

$ws.Cells.Item(141,1).Value = "20100120"
$ws.Cells.Item(141,2).Value = 881

# Normalize the style of the new date cells back to the sheet's default
# (no explicit style) by copying the formatting from the adjacent,
# untouched cell above them.
$ws.Range("A137").Copy()
$ws.Range("A138:A141").PasteSpecial(-4122)
